$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Paragraph 3: "pip install paddlepaddle 2.1.0" -> "pip install paddlepaddle==2.1.0" @ 18pt ---
$para3 = $tr.Paragraphs(3)
$para3.Font.Size = 18

$full = $tr.Text
$idx = $full.IndexOf(" 2.1.0")
$oldRange = $tr.Characters($idx + 1, 6)
$oldRange.Text = "==2.1.0"

# --- Insert the three new paragraphs (requirements) right after paragraph 3 ---
$newText = "`rconda install -c menpo opencv`rconda install -c pytorch pytorch torchvision`rpip install pyyaml"
$para3.InsertAfter($newText)

$para4 = $tr.Paragraphs(4)
$para4.Font.Size = 18
$para5 = $tr.Paragraphs(5)
$para5.Font.Size = 18
$para6 = $tr.Paragraphs(6)
$para6.Font.Size = 18

# --- Split paragraph 4 into its word runs: conda | install -c | menpo | (space) | opencv ---
$f = $para4.Find("conda", 0)
$f.Font.Size = 18
$f = $para4.Find("menpo", 0)
$f.Font.Size = 18
$f = $para4.Find("opencv", 0)
$f.Font.Size = 18

# --- Split paragraph 5 into its word runs: conda | install -c | pytorch | (space) | pytorch | (space) | torchvision ---
$f = $para5.Find("conda", 0)
$f.Font.Size = 18
$after = $f.Start + $f.Length - $para5.Start
$f = $para5.Find("pytorch", $after)
$f.Font.Size = 18
$after = $f.Start + $f.Length - $para5.Start
$f = $para5.Find("pytorch", $after)
$f.Font.Size = 18
$after = $f.Start + $f.Length - $para5.Start
$f = $para5.Find("torchvision", $after)
$f.Font.Size = 18

# --- Split paragraph 6 into its word runs: pip install | pyyaml ---
$f = $para6.Find("pyyaml", 0)
$f.Font.Size = 18

# --- Move the requirements table (graphicFrame "Table 3") down/right ---
$table = $s.Shapes.Item(3)
$table.Left = 130.75464566929134
$table.Top = 414.906062992126
